$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "23.03.2023 10:09 (CET)"
$ws.Range("C4").Value = "https://gitlab.intra.infineon.com/semantic-web-projects/digital-reference/order_management/-/commit/4f24e65057a7257c057b4ab9241d4c02154082ae"
$ws.Range("D4").Value = "54d4d2f24ac8d7a88ec095deed5ffb07daab728d3a812d4d0503aacf7cd69912"
